$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Elastic Collision"
$ws.Range("B4").Value = "collision between pool balls"
$ws.Range("C4").Value = "24/06/2016``"
$ws.Range("D4").Value = "Jack "
$ws.Range("E4").Value = "Research"
$ws.Range("F4").Value = ".svg"
